$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = $null
$ws.Range("B13").Formula = $null
$ws.Range("D13").Formula = $null
$ws.Range("F13").Formula = $null

$ws.Range("A36").Value = "Total"
$ws.Range("B36").Formula = "=SUM(B2:B35)"
$ws.Range("D36").Formula = "=SUM(D2:D35)"
$ws.Range("F36").Formula = "=SUM(F2:F35)"

$ws.Range("F16").Select()
